# Apply "想去人数" (interest count) updates scraped at a later point in time,
# plus two tickets in 演出 that flipped to sold-out ("不可售") status.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1857
$ws.Range("F3").Value = 22
$ws.Range("F4").Value = 25
$ws.Range("F7").Value = 44
$ws.Range("F9").Value = 645
$ws.Range("F10").Value = 70
$ws.Range("F12").Value = 770
$ws.Range("F13").Value = 1471
$ws.Range("F15").Value = 1450
$ws.Range("F16").Value = 29
$ws.Range("F17").Value = 1266
$ws.Range("F18").Value = 303
$ws.Range("F19").Value = 1596
$ws.Range("F20").Value = 777
$ws.Range("F21").Value = 1021
$ws.Range("F22").Value = 329
$ws.Range("F25").Value = 1416
$ws.Range("F26").Value = 104
$ws.Range("F27").Value = 800
$ws.Range("F28").Value = 547
$ws.Range("F29").Value = 1082
$ws.Range("F30").Value = 264414
$ws.Range("F31").Value = 984
$ws.Range("F32").Value = 23
$ws.Range("F33").Value = 555
$ws.Range("F34").Value = 1324
$ws.Range("F35").Value = 1049
$ws.Range("F36").Value = 890
$ws.Range("F37").Value = 1075
$ws.Range("F38").Value = 26
$ws.Range("F39").Value = 52
$ws.Range("F42").Value = 1601
$ws.Range("F43").Value = 98
$ws.Range("F44").Value = 42
$ws.Range("F45").Value = 793

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 134
$ws.Range("F6").Value = 144
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = "不可售"
$ws.Range("F8").Value = 87
$ws.Range("F10").Value = 175
$ws.Range("G10").Value = "不可售"
$ws.Range("F11").Value = 1441
$ws.Range("F14").Value = 2537
$ws.Range("F17").Value = 715
$ws.Range("F23").Value = 437
$ws.Range("F31").Value = 184
$ws.Range("F32").Value = 235
$ws.Range("F34").Value = 144
$ws.Range("F37").Value = 53
$ws.Range("F41").Value = 35
$ws.Range("F43").Value = 42
$ws.Range("F44").Value = 42
$ws.Range("F46").Value = 127

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 79
$ws.Range("F4").Value = 250
$ws.Range("F5").Value = 2795
$ws.Range("F6").Value = 4540
$ws.Range("F7").Value = 123
$ws.Range("F9").Value = 538
$ws.Range("F10").Value = 656
$ws.Range("F11").Value = 434
$ws.Range("F12").Value = 230
$ws.Range("F13").Value = 832
$ws.Range("F14").Value = 212
$ws.Range("F15").Value = 477

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1857
$ws.Range("F3").Value = 250
$ws.Range("F4").Value = 2795
$ws.Range("F5").Value = 4540
$ws.Range("F6").Value = 656
$ws.Range("F8").Value = 230
$ws.Range("F9").Value = 230
$ws.Range("F10").Value = 832
$ws.Range("F11").Value = 832
$ws.Range("F12").Value = 212
$ws.Range("F14").Value = 44
$ws.Range("F16").Value = 1441
$ws.Range("F18").Value = 770
$ws.Range("F19").Value = 2537
$ws.Range("F21").Value = 1471
$ws.Range("F23").Value = 1450
$ws.Range("F24").Value = 1267
$ws.Range("F27").Value = 1596
$ws.Range("F28").Value = 777
$ws.Range("F29").Value = 1021
$ws.Range("F30").Value = 329
$ws.Range("F31").Value = 477
$ws.Range("F32").Value = 477
$ws.Range("F33").Value = 437
$ws.Range("F34").Value = 1416
$ws.Range("F35").Value = 800
$ws.Range("F36").Value = 547
$ws.Range("F37").Value = 1082
$ws.Range("F39").Value = 984
$ws.Range("F40").Value = 24
$ws.Range("F41").Value = 1049
$ws.Range("F42").Value = 890
$ws.Range("F43").Value = 1075
$ws.Range("F47").Value = 1601
$ws.Range("F48").Value = 98
$ws.Range("F49").Value = 793
$ws.Range("F50").Value = 42
